$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 ("Prix Garantie feuchtes Toilettenpapier 2x70 Stück") - all rows below shift up by one
$ws.Rows.Item(3).Delete()

# Update the timestamp column (O) for all remaining data rows (2 through 29) to the new scrape time
$newTimestamp = "2022-08-22 20:58:12"
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}
